$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Cod. Moneda" column (F) held the literal "US$" for every data row
# (rows 2-18). Update it to read "USD" instead.
$ws.Range("F2:F18").Value = "USD"

# Reset the saved selection back to the default top-left cell.
$ws.Range("A1").Select()
